$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block tweaks -----------------------------------------------
# "Valor Mora" (cuenta) total changed
$ws.Range("E11").Value = 898527
# "Cant. Trabajadores" / "Cant. Periodos" counts changed
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 30

# --- Re-style row 45 to become the new "last data row" -----------------
# Row 46 currently carries the distinctive thicker-border "closing" style
# used for the final data row. Copy that formatting onto row 45 (which
# will become the last period row once row 46 is removed below), then
# overwrite row 46 itself with that very same formatting too (so the
# delete below cleanly collapses the table without leaving stray style
# remnants behind).
$ws.Range("B46:J46").Copy() | Out-Null
$ws.Range("B45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 45 reverts to the regular worker identity columns (same worker as
# all the other period rows above it).
$ws.Range("B45").Value = "CC"
$ws.Range("C45").Value = "1047414570"
$ws.Range("D45").Value = "ALBERTO MARIO MORALES DE LUQUE"

# --- Rewrite the 30 period rows (16-45) with the refreshed dataset -----
# Periods now run oldest -> newest (1707 .. 1912) instead of newest ->
# oldest, and the "Valor Mora" / "Salario Basico" figures were updated.
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 781242
$ws.Range("E17").Value = "1708"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242
$ws.Range("E18").Value = "1709"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242
$ws.Range("E19").Value = "1710"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242
$ws.Range("E20").Value = "1711"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242
$ws.Range("E21").Value = "1712"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 781242
$ws.Range("E22").Value = "1801"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 781242
$ws.Range("E23").Value = "1802"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 781242
$ws.Range("E24").Value = "1803"
$ws.Range("F24").Value = 29509
$ws.Range("G24").Value = 781242
$ws.Range("E25").Value = "1804"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 781242
$ws.Range("E26").Value = "1805"
$ws.Range("F26").Value = 29509
$ws.Range("G26").Value = 781242
$ws.Range("E27").Value = "1806"
$ws.Range("F27").Value = 29509
$ws.Range("G27").Value = 781242
$ws.Range("E28").Value = "1807"
$ws.Range("F28").Value = 29509
$ws.Range("G28").Value = 781242
$ws.Range("E29").Value = "1808"
$ws.Range("F29").Value = 29509
$ws.Range("G29").Value = 781242
$ws.Range("E30").Value = "1809"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242
$ws.Range("E31").Value = "1810"
$ws.Range("F31").Value = 31249
$ws.Range("G31").Value = 781242
$ws.Range("E32").Value = "1811"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242
$ws.Range("E33").Value = "1812"
$ws.Range("F33").Value = 31249
$ws.Range("G33").Value = 781242
$ws.Range("E34").Value = "1901"
$ws.Range("F34").Value = 31249
$ws.Range("G34").Value = 781242
$ws.Range("E35").Value = "1902"
$ws.Range("F35").Value = 31249
$ws.Range("G35").Value = 781242
$ws.Range("E36").Value = "1903"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242
$ws.Range("E37").Value = "1904"
$ws.Range("F37").Value = 31249
$ws.Range("G37").Value = 781242
$ws.Range("E38").Value = "1905"
$ws.Range("F38").Value = 31249
$ws.Range("G38").Value = 781242
$ws.Range("E39").Value = "1906"
$ws.Range("F39").Value = 31249
$ws.Range("G39").Value = 781242
$ws.Range("E40").Value = "1907"
$ws.Range("F40").Value = 31249
$ws.Range("G40").Value = 781242
$ws.Range("E41").Value = "1908"
$ws.Range("F41").Value = 31249
$ws.Range("G41").Value = 781242
$ws.Range("E42").Value = "1909"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242
$ws.Range("E43").Value = "1910"
$ws.Range("F43").Value = 31249
$ws.Range("G43").Value = 781242
$ws.Range("E44").Value = "1911"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242
$ws.Range("E45").Value = "1912"
$ws.Range("F45").Value = 16666
$ws.Range("G45").Value = 781242

# --- Drop the now-obsolete 31st row (old worker CAMILO ALBERTO
# FELIZZOLA MERLANO / period 1606) and shift everything below up -------
$ws.Rows(46).Delete()
